$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for rows 2-5 (all 34 columns A:AH) ---
# New sensor dataset ("custom accuracy" / ~1000-row dataset resample) replacing the old readings,
# and timestamps bumped forward to the new sampling window.
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45069.50694444445
$data[0,1] = 5.885
$data[0,2] = 1.942
$data[0,3] = 1.363
$data[0,4] = 7.767
$data[0,5] = 3.737
$data[0,6] = 0.96
$data[0,7] = 6.281
$data[0,8] = 2.168
$data[0,9] = 0.758
$data[0,10] = 1.218
$data[0,11] = 2.359
$data[0,12] = 5.493
$data[0,13] = 0.668
$data[0,14] = 0.884
$data[0,15] = 2.495
$data[0,16] = 1.484
$data[0,17] = 1.256
$data[0,18] = 0
$data[0,19] = 25.797
$data[0,20] = 5.514
$data[0,21] = 3.303
$data[0,22] = 3.938
$data[0,23] = 2.803
$data[0,24] = 0.246
$data[0,25] = 1.46
$data[0,26] = 1.634
$data[0,27] = 0.588
$data[0,28] = 2.048
$data[0,29] = 3.034
$data[0,30] = 0.186
$data[0,31] = 2.547
$data[0,32] = 0.594
$data[0,33] = 1.826
$data[1,0] = 45069.51388888889
$data[1,1] = 20.418
$data[1,2] = 14.511
$data[1,3] = 1.16
$data[1,4] = 42.506
$data[1,5] = 34.077
$data[1,6] = 14.471
$data[1,7] = 51.931
$data[1,8] = 22.718
$data[1,9] = 10.025
$data[1,10] = 15.053
$data[1,11] = 16.759
$data[1,12] = 18.903
$data[1,13] = 4.823
$data[1,14] = 14.271
$data[1,15] = 21.359
$data[1,16] = 12.337
$data[1,17] = 0.754
$data[1,18] = 0.426
$data[1,19] = 221.15
$data[1,20] = 41.825
$data[1,21] = 14.318
$data[1,22] = 28.583
$data[1,23] = 15.256
$data[1,24] = 1.962
$data[1,25] = 25.906
$data[1,26] = 12.225
$data[1,27] = 10.202
$data[1,28] = 12.782
$data[1,29] = 17.869
$data[1,30] = 0.173
$data[1,31] = 45.806
$data[1,32] = 7.61
$data[1,33] = 16.997
$data[2,0] = 45069.52083333334
$data[2,1] = 7.955
$data[2,2] = 5.509
$data[2,3] = 0.555
$data[2,4] = 16.135
$data[2,5] = 12.678
$data[2,6] = 5.236
$data[2,7] = 26.186
$data[2,8] = 8.390000000000001
$data[2,9] = 3.647
$data[2,10] = 5.472
$data[2,11] = 6.285
$data[2,12] = 7.385
$data[2,13] = 1.814
$data[2,14] = 5.098
$data[2,15] = 8.064
$data[2,16] = 4.611
$data[2,17] = 0.459
$data[2,18] = 0.08699999999999999
$data[2,19] = 78.28700000000001
$data[2,20] = 15.933
$data[2,21] = 5.47
$data[2,22] = 10.99
$data[2,23] = 5.844
$data[2,24] = 0.736
$data[2,25] = 12.049
$data[2,26] = 4.584
$data[2,27] = 3.662
$data[2,28] = 4.851
$data[2,29] = 6.771
$data[2,30] = 0.141
$data[2,31] = 23.378
$data[2,32] = 2.76
$data[2,33] = 6.277
$data[3,0] = 45069.52777777778
$data[3,1] = 18.32
$data[3,2] = 13.45
$data[3,3] = 0.82
$data[3,4] = 39.01
$data[3,5] = 31.77
$data[3,6] = 13.66
$data[3,7] = 52.27
$data[3,8] = 21.29
$data[3,9] = 9.42
$data[3,10] = 14.15
$data[3,11] = 15.51
$data[3,12] = 16.95
$data[3,13] = 4.47
$data[3,14] = 13.49
$data[3,15] = 19.86
$data[3,16] = 11.49
$data[3,17] = 0.45
$data[3,18] = 0.42
$data[3,19] = 204.22
$data[3,20] = 38.81
$data[3,21] = 13.03
$data[3,22] = 26.44
$data[3,23] = 13.98
$data[3,24] = 1.82
$data[3,25] = 25.43
$data[3,26] = 11.33
$data[3,27] = 9.67
$data[3,28] = 11.8
$data[3,29] = 16.43
$data[3,30] = 0.12
$data[3,31] = 46.69
$data[3,32] = 7.16
$data[3,33] = 15.88
$ws.Range("A2:AH5").Value = $data

# --- Remove the old 6th data row (dataset now has one fewer sample row; dimension becomes A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- Widen specific data columns from 7 to 8 characters ---
# (ColumnWidth is in characters; the stored OOXML "width" = ColumnWidth + 5/6, so 7.166666666666667 -> stored width 8)
$wideColumns = @(2, 3, 7, 9, 10, 11, 12, 13, 15, 16, 17, 22, 24, 27, 28, 29, 30, 34)
foreach ($colIndex in $wideColumns) {
    $ws.Columns.Item($colIndex).ColumnWidth = 7.166666666666667
}
